$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the "K" column (G) values with the newly computed strike-derived values
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 2
